$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date column header (CO1) and corresponding data (CO2:CO11)
$ws.Range("CO1").Value = "27-nov"

$values = @(12, 8, 8, 12, 10, 16, 15, 9, 14, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 93).Value = $values[$i]
}

# Match formatting of the preceding date column (CN) for header and data cells
$ws.Range("CO1").NumberFormat = $ws.Range("CN1").NumberFormat
$ws.Range("CO2:CO11").NumberFormat = $ws.Range("CN2:CN11").NumberFormat
$ws.Range("CO2:CO11").HorizontalAlignment = $ws.Range("CN2:CN11").HorizontalAlignment

# Update selection as recorded in the saved workbook state
$ws.Range("CT15").Select()
